# Adds extracted data from 10.1016/j.matchar.2023.113025 (hot extrusion study)
# on Ti22 Sc22 Zr22 Nb17 V17 into rows 192-203 of the "Table" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Common values reused across the new rows
# ---------------------------------------------------------------------------
$composition      = "Ti22 Sc22 Zr22 Nb17 V17"
$doi              = "10.1016/j.matchar.2023.113025"
$pointer          = "T4"

$procHomog        = "homogenized at 1173K for 24h"
$procHomogExtr    = "homogenized at 1173K for 24h; hot extruded at 1173K in steel sleeve at 22mm/s to 11.4:1 ratio"
$procHomogExtrAnn = "homogenized at 1173K for 24h; hot extruded at 1173K in steel sleeve at 22mm/s to 11.4:1 ratio; annealed at 1273K for 24h in vacuum"

$nickHomog        = "VAM+H"
$nickHomogExtr    = "VAM+H+HE"
$nickHomogExtrAnn = "VAM+H+HE+A"

$structHomog        = "BCC+BCC+HCP+FCC"
$structHomogExtr    = "BCC+BCC+HCP+HCP+FCC+FCC"
$structHomogExtrAnn = "BCC+BCC+HCP+FCC+FCC"

$expMethod = "EXP"

# ---------------------------------------------------------------------------
# Row 192 - hardness, homogenized
# ---------------------------------------------------------------------------
$ws.Range("B192").Value = $composition
$ws.Range("C192").Value = $structHomog
$ws.Range("D192").Value = $nickHomog
$ws.Range("E192").Value = $procHomog
$ws.Range("F192").Value = "hardness"
$ws.Range("G192").Value = $expMethod
$ws.Range("I192").Value = 298
$ws.Range("P192").Value = 286
$ws.Range("Q192").Value = 13
$ws.Range("J192").Formula = "=P192*9807000"
$ws.Range("K192").Formula = "=Q192*9807000"
$ws.Range("L192").Value = "Pa"
$ws.Range("M192").Value = $pointer
$ws.Range("N192").Value = $doi

# ---------------------------------------------------------------------------
# Row 193 - hardness, homogenized + hot extruded
# ---------------------------------------------------------------------------
$ws.Range("B193").Value = $composition
$ws.Range("C193").Value = $structHomogExtr
$ws.Range("D193").Value = $nickHomogExtr
$ws.Range("E193").Value = $procHomogExtr
$ws.Range("F193").Value = "hardness"
$ws.Range("G193").Value = $expMethod
$ws.Range("I193").Value = 298
$ws.Range("P193").Value = 291
$ws.Range("Q193").Value = 8
$ws.Range("J193").Formula = "=P193*9807000"
$ws.Range("K193").Formula = "=Q193*9807000"
$ws.Range("L193").Value = "Pa"
$ws.Range("M193").Value = $pointer
$ws.Range("N193").Value = $doi

# ---------------------------------------------------------------------------
# Row 194 - hardness, homogenized + hot extruded + annealed
# ---------------------------------------------------------------------------
$ws.Range("B194").Value = $composition
$ws.Range("C194").Value = $structHomogExtrAnn
$ws.Range("D194").Value = $nickHomogExtrAnn
$ws.Range("E194").Value = $procHomogExtrAnn
$ws.Range("F194").Value = "hardness"
$ws.Range("G194").Value = $expMethod
$ws.Range("I194").Value = 298
$ws.Range("P194").Value = 271
$ws.Range("Q194").Value = 7
$ws.Range("J194").Formula = "=P194*9807000"
$ws.Range("K194").Formula = "=Q194*9807000"
$ws.Range("L194").Value = "Pa"
$ws.Range("M194").Value = $pointer
$ws.Range("N194").Value = $doi

# ---------------------------------------------------------------------------
# Row 195 - tensile yield stress, homogenized
# ---------------------------------------------------------------------------
$ws.Range("B195").Value = $composition
$ws.Range("C195").Value = $structHomog
$ws.Range("D195").Value = $nickHomog
$ws.Range("E195").Value = $procHomog
$ws.Range("F195").Value = "tensile yield stress"
$ws.Range("G195").Value = $expMethod
$ws.Range("I195").Value = 298
$ws.Range("J195").Value = 751000000
$ws.Range("K195").Value = 40000000
$ws.Range("L195").Value = "Pa"
$ws.Range("M195").Value = $pointer
$ws.Range("N195").Value = $doi

# ---------------------------------------------------------------------------
# Row 196 - tensile yield stress, homogenized + hot extruded
# ---------------------------------------------------------------------------
$ws.Range("B196").Value = $composition
$ws.Range("C196").Value = $structHomogExtr
$ws.Range("D196").Value = $nickHomogExtr
$ws.Range("E196").Value = $procHomogExtr
$ws.Range("F196").Value = "tensile yield stress"
$ws.Range("G196").Value = $expMethod
$ws.Range("I196").Value = 298
$ws.Range("J196").Value = 936000000
$ws.Range("K196").Value = 35000000
$ws.Range("L196").Value = "Pa"
$ws.Range("M196").Value = $pointer
$ws.Range("N196").Value = $doi

# ---------------------------------------------------------------------------
# Row 197 - tensile yield stress, homogenized + hot extruded + annealed
# ---------------------------------------------------------------------------
$ws.Range("B197").Value = $composition
$ws.Range("C197").Value = $structHomogExtrAnn
$ws.Range("D197").Value = $nickHomogExtrAnn
$ws.Range("E197").Value = $procHomogExtrAnn
$ws.Range("F197").Value = "tensile yield stress"
$ws.Range("G197").Value = $expMethod
$ws.Range("I197").Value = 298
$ws.Range("J197").Value = 740000000
$ws.Range("K197").Value = 45000000
$ws.Range("L197").Value = "Pa"
$ws.Range("M197").Value = $pointer
$ws.Range("N197").Value = $doi

# ---------------------------------------------------------------------------
# Row 198 - UTS, homogenized
# ---------------------------------------------------------------------------
$ws.Range("B198").Value = $composition
$ws.Range("C198").Value = $structHomog
$ws.Range("D198").Value = $nickHomog
$ws.Range("E198").Value = $procHomog
$ws.Range("F198").Value = "UTS"
$ws.Range("G198").Value = $expMethod
$ws.Range("I198").Value = 298
$ws.Range("J198").Value = 848000000
$ws.Range("K198").Value = 50000000
$ws.Range("L198").Value = "Pa"
$ws.Range("M198").Value = $pointer
$ws.Range("N198").Value = $doi

# ---------------------------------------------------------------------------
# Row 199 - UTS, homogenized + hot extruded
# ---------------------------------------------------------------------------
$ws.Range("B199").Value = $composition
$ws.Range("C199").Value = $structHomogExtr
$ws.Range("D199").Value = $nickHomogExtr
$ws.Range("E199").Value = $procHomogExtr
$ws.Range("F199").Value = "UTS"
$ws.Range("G199").Value = $expMethod
$ws.Range("I199").Value = 298
$ws.Range("J199").Value = 939000000
$ws.Range("K199").Value = 40000000
$ws.Range("L199").Value = "Pa"
$ws.Range("M199").Value = $pointer
$ws.Range("N199").Value = $doi

# ---------------------------------------------------------------------------
# Row 200 - UTS, homogenized + hot extruded + annealed
# ---------------------------------------------------------------------------
$ws.Range("B200").Value = $composition
$ws.Range("C200").Value = $structHomogExtrAnn
$ws.Range("D200").Value = $nickHomogExtrAnn
$ws.Range("E200").Value = $procHomogExtrAnn
$ws.Range("F200").Value = "UTS"
$ws.Range("G200").Value = $expMethod
$ws.Range("I200").Value = 298
$ws.Range("J200").Value = 823000000
$ws.Range("K200").Value = 50000000
$ws.Range("L200").Value = "Pa"
$ws.Range("M200").Value = $pointer
$ws.Range("N200").Value = $doi

# ---------------------------------------------------------------------------
# Row 201 - tensile ductility, homogenized
# ---------------------------------------------------------------------------
$ws.Range("B201").Value = $composition
$ws.Range("C201").Value = $structHomog
$ws.Range("D201").Value = $nickHomog
$ws.Range("E201").Value = $procHomog
$ws.Range("F201").Value = "tensile ductility"
$ws.Range("G201").Value = $expMethod
$ws.Range("I201").Value = 298
$ws.Range("J201").Value = 10.2
$ws.Range("K201").Value = 1.2
$ws.Range("L201").Value = "%"
$ws.Range("M201").Value = $pointer
$ws.Range("N201").Value = $doi

# ---------------------------------------------------------------------------
# Row 202 - tensile ductility, homogenized + hot extruded
# ---------------------------------------------------------------------------
$ws.Range("B202").Value = $composition
$ws.Range("C202").Value = $structHomogExtr
$ws.Range("D202").Value = $nickHomogExtr
$ws.Range("E202").Value = $procHomogExtr
$ws.Range("F202").Value = "tensile ductility"
$ws.Range("G202").Value = $expMethod
$ws.Range("I202").Value = 298
$ws.Range("J202").Value = 10.7
$ws.Range("K202").Value = 1.5
$ws.Range("L202").Value = "%"
$ws.Range("M202").Value = $pointer
$ws.Range("N202").Value = $doi

# ---------------------------------------------------------------------------
# Row 203 - tensile ductility, homogenized + hot extruded + annealed
# ---------------------------------------------------------------------------
$ws.Range("B203").Value = $composition
$ws.Range("C203").Value = $structHomogExtrAnn
$ws.Range("D203").Value = $nickHomogExtrAnn
$ws.Range("E203").Value = $procHomogExtrAnn
$ws.Range("F203").Value = "tensile ductility"
$ws.Range("G203").Value = $expMethod
$ws.Range("I203").Value = 298
$ws.Range("J203").Value = 18.5
$ws.Range("K203").Value = 1.7
$ws.Range("L203").Value = "%"
$ws.Range("M203").Value = $pointer
$ws.Range("N203").Value = $doi

# ---------------------------------------------------------------------------
# Update scroll position / selection to match the author's final view
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 173
$win.ScrollColumn = 1
try { $win.TopLeftCell = $ws.Range("A173") } catch {}
$ws.Range("M209").Select()
